# Regenerate the "K" column (column G) values in save_data sheet.
# The column G previously held a "Strike#" derived value; it has been
# regenerated to hold strikeout counts (K) instead, with std/mean recalculated
# upstream and the resulting s_vals written back into the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @{
    2  = 11
    3  = 6
    4  = 3
    5  = 4
    6  = 4
    7  = 6
    8  = 5
    9  = 5
    10 = 5
    11 = 1
    12 = 8
    13 = 12
    14 = 3
    15 = 5
    16 = 5
    17 = 5
    18 = 2
    19 = 2
}

foreach ($row in $newValues.Keys) {
    $ws.Range("G$row").Value = $newValues[$row]
}
